# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" / handoff / handback timestamps
# for the ea18cc3d-a63e-4ec2-8fc9-a03f72b78516 row, across the Overview,
# zh-cn and de-de sheets, reflecting a newly generated handback report.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for the
# ea18cc3d-a63e-4ec2-8fc9-a03f72b78516 row (row 3).
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G3").Value = "2016-08-18 02:42:16"

# zh-cn sheet: "Correspond Handoff Datetime" (H) and
# "Correspond Handback DateTime" (K) for the ea18cc3d row (row 3).
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H3").Value = "2016-08-18 02:42:11"
$zhcn.Range("K3").Value = "2016-08-18 02:42:26"

# de-de sheet: "Correspond Handoff Datetime" (H) and
# "Correspond Handback DateTime" (K) for the ea18cc3d row (row 3).
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H3").Value = "2016-08-18 02:42:16"
$dede.Range("K3").Value = "2016-08-18 02:42:37"
